# Apply Fri Nov  8 13:54:33 UTC 2024 cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.157.04"
$ws.Range("E2").Value = "  +1.88%  "

$ws.Range("D3").Value = "2.944.87"
$ws.Range("E3").Value = "  +4.72%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'203.89"
$ws.Range("E5").Value = "  +9.32%  "

$ws.Range("D6").Value = "'598.53"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.552"
$ws.Range("E8").Value = "  +1.33%  "

$ws.Range("E9").Value = "  +6.05%  "

$ws.Range("D10").Value = "2.931.25"
$ws.Range("E10").Value = "  +4.14%  "

$ws.Range("D11").Value = "'0.438"
$ws.Range("E11").Value = "  +17.68%  "

$ws.Range("D12").Value = "'0.162"
$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("D13").Value = "'4.95"
$ws.Range("E13").Value = "  +2.09%  "

$ws.Range("D14").Value = "3.480.67"
$ws.Range("E14").Value = "  +4.45%  "

$ws.Range("D15").Value = "76.025.04"
$ws.Range("E15").Value = "  +1.80%  "

$ws.Range("D16").Value = "'28.13"
$ws.Range("E16").Value = "  +5.34%  "

$ws.Range("D17").Value = "'0.0000191"
$ws.Range("E17").Value = "  +3.21%  "

$ws.Range("D18").Value = "2.933.15"
$ws.Range("E18").Value = "  +4.04%  "

$ws.Range("D19").Value = "'13.27"
$ws.Range("E19").Value = "  +8.81%  "

$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").Value = "'373.73"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("D22").Value = "'2.31"
$ws.Range("E22").Value = "  +2.83%  "

$ws.Range("E23").Value = "  +6.00%  "

$ws.Range("D24").Value = "'71.82"
$ws.Range("E24").Value = "  +1.55%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +4.60%  "

$ws.Range("D27").Value = "'4.34"
$ws.Range("E27").Value = "  +5.22%  "

$ws.Range("D28").Value = "'9.71"
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("E29").Value = "  +7.40%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  +0.62%  "

$ws.Range("D32").Value = "'505.11"
$ws.Range("E32").Value = "  -0.95%  "

$ws.Range("D33").Value = "'7.82"
$ws.Range("E33").Value = "  +2.91%  "

$ws.Range("E34").Value = "  +3.74%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("D36").Value = "'20.32"
$ws.Range("E36").Value = "  +2.66%  "

$ws.Range("D37").Value = "'163.83"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("E38").Value = "  +28.19%  "

$ws.Range("D39").Value = "'19.65"
$ws.Range("E39").Value = "  +1.57%  "

$ws.Range("D40").Value = "'0.372"
$ws.Range("E40").Value = "  +10.01%  "

$ws.Range("E41").Value = "  -3.73%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'181.89"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "'5.01"
$ws.Range("E44").Value = "  +1.33%  "

$ws.Range("D45").Value = "'1.67"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "'1.21"
$ws.Range("E46").Value = "  +0.42%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'40.15"
$ws.Range("E47").Value = "  +0.46%  "

$ws.Range("E48").Value = "  +2.63%  "

$ws.Range("D49").Value = "'0.584"
$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("D50").Value = "'3.76"
$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("D51").Value = "'22.61"
$ws.Range("E51").Value = "  +9.29%  "

Write-Output "Updated cryptos list"